$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Project was renamed: update ProjectID (A), ProjectName (B), and Key2 (F)
# for all four data rows (2-5).
$ws.Range("A2:A5").Value = "02c54b4d-e85b-4bdf-8bb6-48d07f872f1c"
$ws.Range("B2:B5").Value = "Project Love"
$ws.Range("F2:F5").Value = "Chez Martha"

# Move the selection to F2:F5 (active cell F2)
$null = $ws.Range("F2:F5").Select()
